$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Cells.Item(1, 2).Value = "Período"
$ws.Cells.Item(1, 3).Value = "Valor"
$ws.Cells.Item(1, 4).Value = "Categoria"

# Copy the header style (bold font, border, centered alignment) from A1 onto the new D1 header cell
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# Product names in the order they appear in the sheet
$products = @("Gás-BR", "Gás-NE", "Gás-SE", "LGN-BR", "LGN-NE", "LGN-SE", "Petróleo-BR", "Petróleo-NE", "Petróleo-SE")

# Values for "variação do último ano" (2025-2024) - these were the former column B values
$lastYearValues = @(
    -35.24428059585877,
    -43.81402676005936,
    -38.07047593582452,
    -44.07929141415514,
    -29.18970863480975,
    $null,
    -36.91842092174158,
    -37.84140403772091,
    -26.48295957781719
)

# Values for "variação desde 1997" (2025/1997) - these were the former column C values
$since1997Values = @(
    269.5602310681157,
    -54.68498494151564,
    -98.15290549268987,
    52.62759760343716,
    -86.42868143611607,
    -100,
    153.3741001090112,
    -79.12443946814697,
    -77.64233321010346
)

# First block: rows 2-10, period "2025-2024", category "Variação do último ano"
for ($i = 0; $i -lt $products.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $products[$i]
    $ws.Cells.Item($row, 2).Value = "2025-2024"
    if ($null -eq $lastYearValues[$i]) {
        $ws.Cells.Item($row, 3).Value = ""
    } else {
        $ws.Cells.Item($row, 3).Value = $lastYearValues[$i]
    }
    $ws.Cells.Item($row, 4).Value = "Variação do último ano"
}

# Second block: rows 11-19, period "2025/1997", category "Variação desde 1997"
for ($i = 0; $i -lt $products.Length; $i++) {
    $row = 11 + $i
    $ws.Cells.Item($row, 1).Value = $products[$i]
    $ws.Cells.Item($row, 2).Value = "2025/1997"
    $ws.Cells.Item($row, 3).Value = $since1997Values[$i]
    $ws.Cells.Item($row, 4).Value = "Variação desde 1997"
}

Write-Output "done"
